$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Huaping%Zhu%xref no email%1,           Lin%Wang%xref no email%1,           Chengzhi%Fang%xref no email%1,           Sicong%Peng%xref no email%1,           Lianhong%Zhang%xref no email%1,           Guiping%Chang%xref no email%1,           Shiwen%Xia%xref no email%0,           Wenhao%Zhou%xref no email%0]'
$ws.Range('E3').Value = '[Xiaotong%Wang%NULL%0,                  Zhiqiang%Zhou%NULL%0,                  Jianping%Zhang%NULL%0,                  Fengfeng%Zhu%NULL%0,                  Yongyan%Tang%NULL%0,                  Xinghua%Shen%luckydrx@163.com%0]'
$ws.Range('E4').Value = '[Yan%Chen%NULL%0,                  Hua%Peng%NULL%1,                  Lin%Wang%NULL%2,                  Yin%Zhao%NULL%1,                  Lingkong%Zeng%NULL%1,                  Hui%Gao%NULL%1,                  Yalan%Liu%NULL%1]'
$ws.Range('E5').Value = '[Suliman%Khan%NULL%1,                  Liangyu%Peng%NULL%2,                  Liangyu%Peng%NULL%0,                  Rabeea%Siddique%NULL%1,                  Ghulam%Nabi%NULL%1,                  NULL%Nawsherwan%NULL%1,                  Mengzhou%Xue%NULL%1,                  Jianbo%Liu%NULL%1,                  Guang%Han%NULL%1]'
$ws.Range('E6').Value = '[Na%Li%NULL%0,                  Lefei%Han%NULL%1,                  Min%Peng%2658706528@qq.com%1,                  Yuxia%Lv%NULL%1,                  Yin%Ouyang%NULL%1,                  Kui%Liu%NULL%0,                  Linli%Yue%NULL%1,                  Qiannan%Li%NULL%1,                  Guoqiang%Sun%NULL%2,                  Lin%Chen%NULL%0,                  Lin%Yang%l.yang@polyu.edu.hk%1]'
$ws.Range('E7').Value = '[Nan%Yu%NULL%1,                  Wei%Li%NULL%0,                  Qingling%Kang%NULL%1,                  Zhi%Xiong%NULL%1,                  Shaoshuai%Wang%NULL%1,                  Xingguang%Lin%NULL%1,                  Yanyan%Liu%NULL%1,                  Juan%Xiao%NULL%1,                  Haiyi%Liu%NULL%1,                  Dongrui%Deng%NULL%1,                  Suhua%Chen%NULL%1,                  Wanjiang%Zeng%NULL%1,                  Ling%Feng%NULL%1,                  Jianli%Wu%NULL%1]'
$ws.Range('E8').Value = '[Lingkong%Zeng%xref no email%0,           Shiwen%Xia%xref no email%0,           Wenhao%Yuan%xref no email%0,           Kai%Yan%xref no email%0,           Feifan%Xiao%xref no email%0,           Jianbo%Shao%xref no email%0,           Wenhao%Zhou%xref no email%0]'
$ws.Range('E9').Value = '[Nicoletta%Vendola%NULL%1,                  Viviana%stampini%NULL%2,                  Viviana%stampini%NULL%0,                  Roberta%Amadori%NULL%2,                  Roberta%Amadori%NULL%0,                  Martina%Gerbino%NULL%2,                  Martina%Gerbino%NULL%0,                  Annalisa%Curatolo%NULL%2,                  Annalisa%Curatolo%NULL%0,                  Daniela%surico%NULL%2,                  Daniela%surico%NULL%0]'
$ws.Range('E10').Value = '[Lan%Dong%xref no email%0,           Jinhua%Tian%xref no email%0,           Songming%He%xref no email%0,           Chuchao%Zhu%xref no email%0,           Jian%Wang%xref no email%0,           Chen%Liu%xref no email%0,           Jing%Yang%xref no email%0]'
$ws.Range('E11').Value = '[Siyu%Chen%NULL%1,                  E.%Liao%NULL%2,                  E.%Liao%NULL%0,                  Dongmei%Cao%NULL%1,                  Ying%Gao%NULL%1,                  Guoqiang%Sun%sun821ch@126.com%0,                  Yong%Shao%cqshaoyong@163.com%1]'
$ws.Range('E12').Value = '[Baud%D.%coreGivesNoEmail%1,                 Dubruc%E.%coreGivesNoEmail%1,                 Favre%G.%coreGivesNoEmail%1,                 Gengler%C.%coreGivesNoEmail%1,                 Greub%G.%coreGivesNoEmail%1,                 Jaton%K.%coreGivesNoEmail%1,                 Pomar%L.%coreGivesNoEmail%1]'
$ws.Range('E13').Value = '[AZAP%ALPAY%coreGivesNoEmail%0,                 Kalafat%Erkan%coreGivesNoEmail%1,                 KO\u00c7%ASLI%coreGivesNoEmail%1,                 Ozisik%S.%coreGivesNoEmail%1,                 UZUN%\u00c7A\u011eLAR%coreGivesNoEmail%1,                 VARLI%BULUT%coreGivesNoEmail%1,                 Yaprak%E.%coreGivesNoEmail%1,                 \u00c7INAR%G\u00dcLE%coreGivesNoEmail%1]'
$ws.Range('E14').Value = '[Zhoujie%Peng%NULL%0,                  Jianhui%Wang%NULL%1,                  Yunbo%Mo%NULL%1,                  Wei%Duan%NULL%0,                  Guangjun%Xiang%NULL%1,                  Ming%Yi%NULL%1,                  Lei%Bao%NULL%1,                  Yuan%Shi%NULL%1]'
$ws.Range('E15').Value = '[Noelle%Breslin%nb2565@cumc.columbia.edu%0,                  Caitlin%Baptiste%NULL%1,                  Cynthia%Gyamfi-Bannerman%NULL%1,                  Russell%Miller%NULL%1,                  Rebecca%Martinez%NULL%1,                  Kyra%Bernstein%NULL%1,                  Laurence%Ring%NULL%1,                  Ruth%Landau%NULL%1,                  Stephanie%Purisch%NULL%1,                  Alexander M.%Friedman%NULL%1,                  Karin%Fuchs%NULL%1,                  Desmond%Sutton%NULL%1,                  Maria%Andrikopoulou%NULL%1,                  Devon%Rupley%NULL%1,                  Jean-Ju%Sheen%NULL%1,                  Janice%Aubey%NULL%1,                  Noelia%Zork%NULL%1,                  Leslie%Moroz%NULL%1,                  Mirella%Mourad%NULL%1,                  Ronald%Wapner%NULL%1,                  Lynn L.%Simpson%NULL%1,                  Mary E.%D’Alton%NULL%1,                  Dena%Goffman%NULL%1]'
$ws.Range('E16').Value = '[Xiali%Xiong%NULL%0,                  Hong%Wei%NULL%2,                  Hong%Wei%NULL%0,                  Zhihong%Zhang%NULL%1,                  Jing%Chang%NULL%1,                  Xiaopeng%Ma%NULL%1,                  Xiang%Gao%NULL%1,                  Qiang%Chen%qiangchen2015@sina.com%1,                  Qiumei%Pang%pqm19650930@sina.com%2,                  Qiumei%Pang%pqm19650930@sina.com%0]'
$ws.Range('E17').Value = '[Mohammad%Khassawneh%xref no email%1,           Wasim%Khasawneh%xref no email%1,           Laila Al%Zaghal%xref no email%1,           Wail%Hayajneh%xref no email%1,           Fadel%Abdelal%xref no email%1]'
$ws.Range('E18').Value = '[William T.%Schnettler%william_schnettler@trihealth.com%0,                  Yousef%Al Ahwel%NULL%1,                  Anju%Suhag%NULL%1]'
$ws.Range('E19').Value = '[Wei%Liu%NULL%0,                  Jing%Wang%NULL%0,                  Wenbin%Li%NULL%1,                  Zhaoxian%Zhou%NULL%1,                  Siying%Liu%NULL%1,                  Zhihui%Rong%rongzhihui53@163.com%1]'
